# Updates cryptos list cell values to match the latest scraped market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.455.62'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '2.552.71'
$ws.Range('E3').Value = '  +1.09%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = '2.549.75'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.346'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').Value = '2.947.02'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').Value = '68.414.31'
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('B18').Value = 'Binance-PegBSC-USD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +92.83%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '2.538.64'
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.44%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '369.33'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.96%  '
$ws.Range('D29').Value = '2.672.00'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '535.96'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.64'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.44'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.65'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.17'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.349'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.996'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '149.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').Value = '0.0₆0280'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.557'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('E51').Value = '  +2.18%  '
